# Add "Terman et al., 1992" as a new row in the VEGFA165_VEGFR2 sheet.
# This inserts a new row just below the header/group rows (new row 4),
# pushing the existing reference rows (old rows 4-10) down by one, while
# leaving the raw-data columns G:J (which are anchored to their own row
# numbers and referenced by fixed formulas) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEGFA165_VEGFR2")

$xlPasteFormats = -4122

# --- 1. Shift rows 4..10 (columns A:D only) down to rows 5..11. ---
# Work from the bottom up so we never clobber a row before reading it.
for ($r = 9; $r -ge 4; $r--) {
    $rNext = $r + 1
    $srcRow = $ws.Range("A" + $r + ":D" + $r)
    $dstRow = $ws.Range("A" + $rNext + ":D" + $rNext)

    # Copy the formatting (border/number-format styles) down one row.
    $srcRow.Copy()
    $dstRow.PasteSpecial($xlPasteFormats)

    # Copy the actual cell contents (values / formulas) down one row.
    $ws.Range("A" + $rNext).Value = $ws.Range("A" + $r).Value()
    $ws.Range("B" + $rNext).Value = $ws.Range("B" + $r).Value()
    $ws.Range("C" + $rNext).Formula = $ws.Range("C" + $r).Formula()
    $ws.Range("D" + $rNext).Formula = $ws.Range("D" + $r).Formula()
}

# --- 2. Populate the new row 4 with the Terman et al., 1992 entry. ---
# Match the formatting used by the other rows in this first reference
# group (rows 2/3), then fill in the new data.
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial($xlPasteFormats)

$ws.Range("A4").Value = "Terman et al., 1992"
$ws.Range("B4").Value = "Radioligand"
$ws.Range("C4").Value = 75
$ws.Range("D4").ClearContents()

$excel.CutCopyMode = 0

# --- 3. Move the cell comments down one row to match (D6->D7, D7->D8). ---
# Process bottom-up so we never overwrite a comment before it is read.
$commentD7 = $ws.Range("D7").Comment
$textD7 = $commentD7.Text()
$commentD7.Delete()
$ws.Range("D8").AddComment($textD7)

$commentD6 = $ws.Range("D6").Comment
$textD6 = $commentD6.Text()
$commentD6.Delete()
$ws.Range("D7").AddComment($textD6)

# --- 4. Make VEGFA165_VEGFR2 the active/selected sheet. ---
$ws.Activate()
$ws.Range("A5").Select()
